$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 830.875
$ws.Range("J2").Value = 912.25
$ws.Range("L2").Value = 912.25
$ws.Range("N2").Value = -1138.25
$ws.Range("H9").Value = 328.57144
$ws.Range("I9").Value = 325
$ws.Range("K9").Value = 325
$ws.Range("M9").Value = -156
$ws.Range("H20").Value = 2221
$ws.Range("I20").Value = 2221
$ws.Range("K20").Value = 2221
$ws.Range("M20").Value = -1991
$ws.Range("H35").Value = 2221
$ws.Range("I35").Value = 2221
$ws.Range("K35").Value = 2221
$ws.Range("M35").Value = -1842
$ws.Range("H51").Value = 7280.857
$ws.Range("I51").Value = 3000
$ws.Range("J51").Value = 7494.9
$ws.Range("K51").Value = 3000
$ws.Range("L51").Value = 7494.9
$ws.Range("N51").Value = -8462.9
$ws.Range("M51").Value = -2516
$ws.Range("H55").Value = 2427
$ws.Range("I55").Value = 170.33333
$ws.Range("J55").Value = 3555.3333
$ws.Range("K55").Value = 170.33333
$ws.Range("L55").Value = 3555.3333
$ws.Range("M55").Value = 43.66667000000001
$ws.Range("N55").Value = -3983.3333
$ws.Range("H75").Value = 44222
$ws.Range("J75").Value = 44222
$ws.Range("L75").Value = 44222
$ws.Range("N75").Value = -46094
$ws.Range("H78").Value = 44222
$ws.Range("J78").Value = 44222
$ws.Range("L78").Value = 132666
$ws.Range("N78").Value = -142026
$ws.Range("H98").Value = 1444.9
$ws.Range("I98").Value = 1205.1578
$ws.Range("K98").Value = 1205.1578
$ws.Range("M98").Value = 292.8422
$ws.Range("H113").Value = 5998.7085
$ws.Range("I113").Value = 6932.3335
$ws.Range("K113").Value = 6932.3335
$ws.Range("M113").Value = -3678.3335
$ws.Range("H122").Value = 1444.9
$ws.Range("I122").Value = 1205.1578
$ws.Range("K122").Value = 3615.4734
$ws.Range("M122").Value = -1165.4734
$ws.Range("H132").Value = 16394762
$ws.Range("I132").Value = 18519834
$ws.Range("J132").Value = 1349.1428
$ws.Range("K132").Value = 55559502
$ws.Range("L132").Value = 4047.4284
$ws.Range("M132").Value = -55556972
$ws.Range("N132").Value = -9107.428400000001
$ws.Range("H137").Value = 3000.1428
$ws.Range("I137").Value = 2712.375
$ws.Range("J137").Value = 3276.4
$ws.Range("K137").Value = 8137.125
$ws.Range("L137").Value = 9829.200000000001
$ws.Range("M137").Value = -5587.125
$ws.Range("N137").Value = -14929.2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H27").Value = 12040.2
$ws.Range("J27").Value = 12040.2
$ws.Range("L27").Value = 12040.2
$ws.Range("N27").Value = -12408.2
$ws.Range("H32").Value = 5190.1875
$ws.Range("I32").Value = 3258.7917
$ws.Range("K32").Value = 3258.7917
$ws.Range("M32").Value = -2971.7917
$ws.Range("H94").Value = 20164
$ws.Range("J94").Value = 20164
$ws.Range("L94").Value = 20164
$ws.Range("N94").Value = -21966
$ws.Range("H132").Value = 6233.8965
$ws.Range("I132").Value = 6339.522
$ws.Range("K132").Value = 19018.566
$ws.Range("M132").Value = -16488.566
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H57").Value = 24995
$ws.Range("J57").Value = 24995
$ws.Range("L57").Value = 24995
$ws.Range("N57").Value = -26435
$ws.Range("H95").Value = 14000
$ws.Range("J95").Value = 14000
$ws.Range("L95").Value = 14000
$ws.Range("N95").Value = -19492
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()
$ws.Range("H136").Value = 24995
$ws.Range("J136").Value = 24995
$ws.Range("L136").Value = 24995
$ws.Range("N136").Value = -35195
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8794.058999999999
$ws.Range("I31").Value = 12633.6
$ws.Range("K31").Value = 12633.6
$ws.Range("M31").Value = -12338.6
$ws.Range("H34").Value = 8794.058999999999
$ws.Range("I34").Value = 12633.6
$ws.Range("K34").Value = 12633.6
$ws.Range("M34").Value = -12431.6
$ws.Range("H58").Value = 1881
$ws.Range("J58").Value = 2726
$ws.Range("L58").Value = 2726
$ws.Range("N58").Value = -3132
$ws.Range("H88").Value = 35915.668
$ws.Range("I88").Value = 10000
$ws.Range("J88").Value = 41098.8
$ws.Range("K88").Value = 10000
$ws.Range("L88").Value = 41098.8
$ws.Range("N88").Value = -41910.8
$ws.Range("M88").Value = -9594
$ws.Range("H91").Value = 35915.668
$ws.Range("I91").Value = 10000
$ws.Range("J91").Value = 41098.8
$ws.Range("K91").Value = 10000
$ws.Range("L91").Value = 41098.8
$ws.Range("N91").Value = -43906.8
$ws.Range("M91").Value = -8596
$ws.Range("H122").Value = 3091.4666
$ws.Range("I122").Value = 2055.2856
$ws.Range("J122").Value = 3998.125
$ws.Range("K122").Value = 6165.8568
$ws.Range("L122").Value = 11994.375
$ws.Range("M122").Value = -3715.8568
$ws.Range("N122").Value = -16894.375
$ws.Range("H132").Value = 20587.264
$ws.Range("I132").Value = 22862.191
$ws.Range("K132").Value = 68586.573
$ws.Range("M132").Value = -66056.573
$ws.Range("H136").Value = 1881
$ws.Range("J136").Value = 2726
$ws.Range("L136").Value = 8178
$ws.Range("N136").Value = -13278
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 17863426
$ws.Range("I56").Value = 17863426
$ws.Range("K56").Value = 17863426
$ws.Range("M56").Value = -17862896
$ws.Range("H132").Value = 2131.5217
$ws.Range("I132").Value = 829.5714
$ws.Range("J132").Value = 2701.125
$ws.Range("K132").Value = 7466.1426
$ws.Range("L132").Value = 24310.125
$ws.Range("M132").Value = -4936.1426
$ws.Range("N132").Value = -29370.125
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1751715.5
$ws.Range("I80").Value = 4083324
$ws.Range("J80").Value = 3009.125
$ws.Range("K80").Value = 4083324
$ws.Range("L80").Value = 3009.125
$ws.Range("M80").Value = -4082326
$ws.Range("N80").Value = -5005.125
$ws.Range("H83").Value = 1751715.5
$ws.Range("I83").Value = 4083324
$ws.Range("J83").Value = 3009.125
$ws.Range("K83").Value = 20416620
$ws.Range("L83").Value = 15045.625
$ws.Range("M83").Value = -20411628
$ws.Range("N83").Value = -25029.625
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6190.231
$ws.Range("I7").Value = 4588.636
$ws.Range("J7").Value = 14999
$ws.Range("K7").Value = 4588.636
$ws.Range("L7").Value = 14999
$ws.Range("M7").Value = -4476.636
$ws.Range("N7").Value = -15223
$ws.Range("H40").Value = 7232.625
$ws.Range("I40").Value = 6123.143
$ws.Range("J40").Value = 14999
$ws.Range("K40").Value = 6123.143
$ws.Range("L40").Value = 14999
$ws.Range("M40").Value = -5987.143
$ws.Range("N40").Value = -15271
$ws.Range("H122").Value = 5926
$ws.Range("I122").Value = 3862.4443
$ws.Range("J122").Value = 8247.5
$ws.Range("K122").Value = 11587.3329
$ws.Range("L122").Value = 24742.5
$ws.Range("M122").Value = -9137.332900000001
$ws.Range("N122").Value = -29642.5
$ws.Range("H126").Value = 6190.231
$ws.Range("I126").Value = 4588.636
$ws.Range("J126").Value = 14999
$ws.Range("K126").Value = 13765.908
$ws.Range("L126").Value = 44997
$ws.Range("M126").Value = -11295.908
$ws.Range("N126").Value = -49937
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 3658.5833
$ws.Range("I122").Value = 2442
$ws.Range("J122").Value = 6613.143
$ws.Range("K122").Value = 7326
$ws.Range("L122").Value = 19839.429
$ws.Range("M122").Value = -4876
$ws.Range("N122").Value = -24739.429
$ws.Range("H132").Value = 21051278
$ws.Range("I132").Value = 26324750
$ws.Range("K132").Value = 78974250
$ws.Range("M132").Value = -78971720
$ws.Range("H136").Value = 4230.4053
$ws.Range("I136").Value = 4503.2666
$ws.Range("K136").Value = 13509.7998
$ws.Range("M136").Value = -10959.7998
